$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AP3").Value = 2.09
$ws.Range("AQ3").Value = 1.81
$ws.Range("G4").Value = 2.38
$ws.Range("I4").Value = 3.2
$ws.Range("J4").Value = 3.25
$ws.Range("L4").Value = 4
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48
$ws.Range("S4").Value = 5.5
$ws.Range("T4").Value = 1.14
$ws.Range("W4").Value = 2.2
$ws.Range("X4").Value = 1.62
$ws.Range("Z4").Value = 10
$ws.Range("AA4").Value = 11
$ws.Range("AB4").Value = 23
$ws.Range("AM4").Value = 34
$ws.Range("AP4").Value = 2
$ws.Range("AQ4").Value = 1.85
$ws.Range("G5").Value = 1.91
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 2.75
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("Z5").Value = 8
$ws.Range("AB5").Value = 17
$ws.Range("AC5").Value = 21
$ws.Range("AJ5").Value = 8.5
$ws.Range("AK5").Value = 19
$ws.Range("I7").Value = 4.85
$ws.Range("J7").Value = 2.37
$ws.Range("K7").Value = 1.98
$ws.Range("Y7").Value = 5.4
$ws.Range("AC7").Value = 16
$ws.Range("AF7").Value = 6.6
$ws.Range("AG7").Value = 19.5
$ws.Range("AJ7").Value = 11
$ws.Range("AK7").Value = 27
$ws.Range("G8").Value = 1.75
$ws.Range("H8").Value = 3.3
$ws.Range("J8").Value = 2.5
$ws.Range("M8").Value = 1.11
$ws.Range("N8").Value = 6.5
$ws.Range("Q8").Value = 2.6
$ws.Range("R8").Value = 1.48
$ws.Range("U8").Value = 1.57
$ws.Range("V8").Value = 2.25
$ws.Range("AB8").Value = 13
$ws.Range("AE8").Value = 6.5
$ws.Range("AF8").Value = 7
$ws.Range("AJ8").Value = 9.5
$ws.Range("AL8").Value = 19
$ws.Range("AP8").Value = 2.03
$ws.Range("AQ8").Value = 1.83
$ws.Range("M10").Value = 1.11
$ws.Range("N10").Value = 6.5
$ws.Range("G16").Value = 1.7
$ws.Range("H16").Value = 3.75
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 2.3
$ws.Range("L16").Value = 5
$ws.Range("Q16").Value = 1.85
$ws.Range("R16").Value = 2
$ws.Range("S16").Value = 3
$ws.Range("T16").Value = 1.36
$ws.Range("W16").Value = 1.8
$ws.Range("X16").Value = 1.95
$ws.Range("AB16").Value = 13
$ws.Range("AI16").Value = 251
$ws.Range("AK16").Value = 26
$ws.Range("AN16").Value = 41
$ws.Range("O23").Value = 1.36
$ws.Range("P23").Value = 3
$ws.Range("Q23").Value = 2.25
$ws.Range("R23").Value = 1.62
$ws.Range("S23").Value = 4
$ws.Range("T23").Value = 1.22
$ws.Range("Q27").Value = 1.6
$ws.Range("R27").Value = 2.3
$ws.Range("AR27").Value = 2
$ws.Range("AS27").Value = 1.85

$wb.Save()
